# Replace Product, IT, and Finance templates with correct industry-specific content
# This script updates the "Instructions & User Guide" and "KPI Dashboard" sheets
# of the Product KPI Dashboard template, replacing leftover
# "Artificial Intelligence / Machine Learning" boilerplate text with the
# correct Product Development wording.

$wb = $excel.ActiveWorkbook

# --- Sheet: Instructions & User Guide ---
$ws1 = $wb.Worksheets.Item("Instructions & User Guide")

$ws1.Range("A1").Value = "Product Development KPI Dashboard - User Guide & Instructions"
$ws1.Range("B23").Value = "Availability and reliability of Product systems"

# --- Sheet: KPI Dashboard ---
$ws2 = $wb.Worksheets.Item("KPI Dashboard")

$ws2.Range("A1").Value = "PRODUCT DEVELOPMENT - KPI DASHBOARD"
$ws2.Range("A2").Value = "Project: Product Development Implementation"

# Owner column: "ML Engineers" -> "Product Engineers"
$ws2.Range("I10").Value = "Product Engineers"
$ws2.Range("I16").Value = "Product Engineers"
$ws2.Range("I22").Value = "Product Engineers"

# Notes column: update "Critical KPI for Artificial Intelligence and Machine
# Learning success" -> "Critical KPI for Product Development success" for
# every KPI row (rows 8 through 22)
for ($row = 8; $row -le 22; $row++) {
    $cell = $ws2.Cells.Item($row, 11)
    $cell.Value = "Critical KPI for Product Development success"
}
